$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.023264501399539
$ws.Cells.Item(2, 4).Value = 1.024894917242775
$ws.Cells.Item(2, 5).Value = 1.023894699122305
$ws.Cells.Item(2, 6).Value = 1.021726742568586
$ws.Cells.Item(2, 9).Value = 1.029585543477908
$ws.Cells.Item(2, 10).Value = 1.028445695229218
$ws.Cells.Item(2, 11).Value = 1.02772207931628
$ws.Cells.Item(2, 12).Value = 1.026724795664891
$ws.Cells.Item(2, 13).Value = 1.024563220987268
$ws.Cells.Item(2, 14).Value = 1.013494321952103
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.024515217912774
$ws.Cells.Item(3, 4).Value = 1.025991700188154
$ws.Cells.Item(3, 5).Value = 1.024963262948683
$ws.Cells.Item(3, 6).Value = 1.023619608069187
$ws.Cells.Item(3, 9).Value = 1.029822798988809
$ws.Cells.Item(3, 10).Value = 1.029333254455474
$ws.Cells.Item(3, 11).Value = 1.02862526326964
$ws.Cells.Item(3, 12).Value = 1.02759962297843
$ws.Cells.Item(3, 13).Value = 1.026259631484754
$ws.Cells.Item(3, 14).Value = 1.013790809329986
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.025322858030693
$ws.Cells.Item(4, 4).Value = 1.026700174599818
$ws.Cells.Item(4, 5).Value = 1.025653643112468
$ws.Cells.Item(4, 6).Value = 1.024842413008299
$ws.Cells.Item(4, 9).Value = 1.029973572292023
$ws.Cells.Item(4, 10).Value = 1.029905464905022
$ws.Cells.Item(4, 11).Value = 1.0292078873609
$ws.Cells.Item(4, 12).Value = 1.028164057382676
$ws.Cells.Item(4, 13).Value = 1.027354925478569
$ws.Cells.Item(4, 14).Value = 1.013981864607071
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.025661998382019
$ws.Cells.Item(5, 4).Value = 1.026997729594704
$ws.Cells.Item(5, 5).Value = 1.02594363062108
$ws.Cells.Item(5, 6).Value = 1.025356013519523
$ws.Cells.Item(5, 9).Value = 1.030036300843509
$ws.Cells.Item(5, 10).Value = 1.030145523212541
$ws.Cells.Item(5, 11).Value = 1.029452396052058
$ws.Cells.Item(5, 12).Value = 1.028400956970106
$ws.Cells.Item(5, 13).Value = 1.027814825207656
$ws.Cells.Item(5, 14).Value = 1.014061995848141
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.025718918727311
$ws.Cells.Item(6, 4).Value = 1.027047673630082
$ws.Cells.Item(6, 5).Value = 1.025992306341253
$ws.Cells.Item(6, 6).Value = 1.025442222327876
$ws.Cells.Item(6, 9).Value = 1.030046794773554
$ws.Cells.Item(6, 10).Value = 1.030185800906097
$ws.Cells.Item(6, 11).Value = 1.029493425231196
$ws.Cells.Item(6, 12).Value = 1.028440710760087
$ws.Cells.Item(6, 13).Value = 1.027892011794285
$ws.Cells.Item(6, 14).Value = 1.014075439222429
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.025327391172882
$ws.Cells.Item(7, 4).Value = 1.026704151669813
$ws.Cells.Item(7, 5).Value = 1.025657518911102
$ws.Cells.Item(7, 6).Value = 1.02484927758068
$ws.Cells.Item(7, 9).Value = 1.029974413052513
$ws.Cells.Item(7, 10).Value = 1.029908674531949
$ws.Cells.Item(7, 11).Value = 1.029211156169415
$ws.Cells.Item(7, 12).Value = 1.028167224369545
$ws.Cells.Item(7, 13).Value = 1.027361072874294
$ws.Cells.Item(7, 14).Value = 1.013982936063471
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.023687533134421
$ws.Cells.Item(8, 4).Value = 1.025265834578697
$ws.Cells.Item(8, 5).Value = 1.024256045266322
$ws.Cells.Item(8, 6).Value = 1.022366868735116
$ws.Cells.Item(8, 9).Value = 1.029666294495696
$ws.Cells.Item(8, 10).Value = 1.028746087033459
$ws.Cells.Item(8, 11).Value = 1.028027688661279
$ws.Cells.Item(8, 12).Value = 1.027020788708765
$ws.Cells.Item(8, 13).Value = 1.025137034422796
$ws.Cells.Item(8, 14).Value = 1.013594685929047
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02078496913178
$ws.Cells.Item(9, 4).Value = 1.022721838241662
$ws.Cells.Item(9, 5).Value = 1.021778250528735
$ws.Cells.Item(9, 6).Value = 1.017976580469046
$ws.Cells.Item(9, 9).Value = 1.02910226330172
$ws.Cells.Item(9, 10).Value = 1.02668120965962
$ws.Cells.Item(9, 11).Value = 1.025928348325299
$ws.Cells.Item(9, 12).Value = 1.024987930290397
$ws.Cells.Item(9, 13).Value = 1.021199096121808
$ws.Cells.Item(9, 14).Value = 1.012904421137528
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.018840906430333
$ws.Cells.Item(10, 4).Value = 1.021019208352412
$ws.Cells.Item(10, 5).Value = 1.020120631558
$ws.Cells.Item(10, 6).Value = 1.015038135576485
$ws.Cells.Item(10, 9).Value = 1.028711993065817
$ws.Cells.Item(10, 10).Value = 1.025293460049166
$ws.Cells.Item(10, 11).Value = 1.024519199521324
$ws.Cells.Item(10, 12).Value = 1.023623941454287
$ws.Cells.Item(10, 13).Value = 1.018560339152245
$ws.Cells.Item(10, 14).Value = 1.012440054604932
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.017996891673836
$ws.Cells.Item(11, 4).Value = 1.020280321616956
$ws.Cells.Item(11, 5).Value = 1.019401446711802
$ws.Cells.Item(11, 6).Value = 1.0137628181422
$ws.Cells.Item(11, 9).Value = 1.028539605060288
$ws.Cells.Item(11, 10).Value = 1.024689846586931
$ws.Cells.Item(11, 11).Value = 1.023906695924337
$ws.Cells.Item(11, 12).Value = 1.023031194638647
$ws.Cells.Item(11, 13).Value = 1.017414364801076
$ws.Cells.Item(11, 14).Value = 1.012237967560814
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.017683046499482
$ws.Cells.Item(12, 4).Value = 1.020005615320948
$ws.Cells.Item(12, 5).Value = 1.019134090733851
$ws.Cells.Item(12, 6).Value = 1.013288648566199
$ws.Cells.Item(12, 9).Value = 1.028475060337174
$ws.Cells.Item(12, 10).Value = 1.024465225897287
$ws.Cells.Item(12, 11).Value = 1.023678829699267
$ws.Cells.Item(12, 12).Value = 1.022810697614525
$ws.Cells.Item(12, 13).Value = 1.01698817709538
$ws.Cells.Item(12, 14).Value = 1.012162749657158
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.017750382869942
$ws.Cells.Item(13, 4).Value = 1.020064552186347
$ws.Cells.Item(13, 5).Value = 1.019191449470488
$ws.Cells.Item(13, 6).Value = 1.01339038062758
$ws.Cells.Item(13, 9).Value = 1.028488928600426
$ws.Cells.Item(13, 10).Value = 1.024513426509368
$ws.Cells.Item(13, 11).Value = 1.023727723922842
$ws.Cells.Item(13, 12).Value = 1.022858009727982
$ws.Cells.Item(13, 13).Value = 1.017079619655879
$ws.Cells.Item(13, 14).Value = 1.012178891136602
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.017970956118453
$ws.Cells.Item(14, 4).Value = 1.020257619455287
$ws.Cells.Item(14, 5).Value = 1.019379351471764
$ws.Cells.Item(14, 6).Value = 1.013723632614848
$ws.Cells.Item(14, 9).Value = 1.028534280230079
$ws.Cells.Item(14, 10).Value = 1.024671287798816
$ws.Cells.Item(14, 11).Value = 1.023887867702379
$ws.Cells.Item(14, 12).Value = 1.023012974930248
$ws.Cells.Item(14, 13).Value = 1.017379146722815
$ws.Cells.Item(14, 14).Value = 1.012231753171459
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.018106813257575
$ws.Cells.Item(15, 4).Value = 1.020376541180415
$ws.Cells.Item(15, 5).Value = 1.01949509497631
$ws.Cells.Item(15, 6).Value = 1.013928898726074
$ws.Cells.Item(15, 9).Value = 1.028562154954595
$ws.Cells.Item(15, 10).Value = 1.024768496640064
$ws.Cells.Item(15, 11).Value = 1.023986490374308
$ws.Cells.Item(15, 12).Value = 1.023108410979773
$ws.Cells.Item(15, 13).Value = 1.017563625646734
$ws.Cells.Item(15, 14).Value = 1.012264302792218
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.018896873519813
$ws.Cells.Item(16, 4).Value = 1.021068210925118
$ws.Cells.Item(16, 5).Value = 1.020168331072939
$ws.Cells.Item(16, 6).Value = 1.015122710608704
$ws.Cells.Item(16, 9).Value = 1.02872336215611
$ws.Cells.Item(16, 10).Value = 1.025333462407118
$ws.Cells.Item(16, 11).Value = 1.024559799823119
$ws.Cells.Item(16, 12).Value = 1.023663234847182
$ws.Cells.Item(16, 13).Value = 1.018636321331301
$ws.Cells.Item(16, 14).Value = 1.012453444976619
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.019391857735051
$ws.Cells.Item(17, 4).Value = 1.021501635394094
$ws.Cells.Item(17, 5).Value = 1.020590249554736
$ws.Cells.Item(17, 6).Value = 1.015870755503977
$ws.Cells.Item(17, 9).Value = 1.02882357222444
$ws.Cells.Item(17, 10).Value = 1.025687121921665
$ws.Cells.Item(17, 11).Value = 1.024918793868426
$ws.Cells.Item(17, 12).Value = 1.024010687705817
$ws.Cells.Item(17, 13).Value = 1.019308281016089
$ws.Cells.Item(17, 14).Value = 1.012571816505837
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.019680359635778
$ws.Cells.Item(18, 4).Value = 1.02175428696971
$ws.Cells.Item(18, 5).Value = 1.020836210302142
$ws.Cells.Item(18, 6).Value = 1.016306793789644
$ws.Cells.Item(18, 9).Value = 1.028881695272093
$ws.Cells.Item(18, 10).Value = 1.025893144510108
$ws.Cells.Item(18, 11).Value = 1.025127964250335
$ws.Cells.Item(18, 12).Value = 1.024213145801276
$ws.Cells.Item(18, 13).Value = 1.019699899079208
$ws.Cells.Item(18, 14).Value = 1.012640762955883
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.019778695224473
$ws.Cells.Item(19, 4).Value = 1.02184040804759
$ws.Cells.Item(19, 5).Value = 1.020920053434108
$ws.Cells.Item(19, 6).Value = 1.016455423961593
$ws.Cells.Item(19, 9).Value = 1.028901458177205
$ws.Cells.Item(19, 10).Value = 1.025963348759171
$ws.Cells.Item(19, 11).Value = 1.025199247958927
$ws.Cells.Item(19, 12).Value = 1.024282144091747
$ws.Cells.Item(19, 13).Value = 1.019833376094849
$ws.Cells.Item(19, 14).Value = 1.012664255395057
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.01933877279916
$ws.Cells.Item(20, 4).Value = 1.021455149387311
$ws.Cells.Item(20, 5).Value = 1.020544995916968
$ws.Cells.Item(20, 6).Value = 1.015790526827016
$ws.Cells.Item(20, 9).Value = 1.028812854549187
$ws.Cells.Item(20, 10).Value = 1.025649204637655
$ws.Cells.Item(20, 11).Value = 1.024880300488655
$ws.Cells.Item(20, 12).Value = 1.023973430567931
$ws.Cells.Item(20, 13).Value = 1.019236219790713
$ws.Cells.Item(20, 14).Value = 1.012559126474052
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.017906012247366
$ws.Cells.Item(21, 4).Value = 1.02020077292186
$ws.Cells.Item(21, 5).Value = 1.019324025090514
$ws.Cells.Item(21, 6).Value = 1.013625511039302
$ws.Cells.Item(21, 9).Value = 1.028520939455896
$ws.Cells.Item(21, 10).Value = 1.024624812989645
$ws.Cells.Item(21, 11).Value = 1.023840719191246
$ws.Cells.Item(21, 12).Value = 1.022967350532685
$ws.Cells.Item(21, 13).Value = 1.017290957990563
$ws.Cells.Item(21, 14).Value = 1.012216190877266
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.017003205408162
$ws.Cells.Item(22, 4).Value = 1.019410642731504
$ws.Cells.Item(22, 5).Value = 1.018555084579662
$ws.Cells.Item(22, 6).Value = 1.012261610575769
$ws.Cells.Item(22, 9).Value = 1.028334436953043
$ws.Cells.Item(22, 10).Value = 1.02397835247459
$ws.Cells.Item(22, 11).Value = 1.023185036232645
$ws.Cells.Item(22, 12).Value = 1.022332908631289
$ws.Cells.Item(22, 13).Value = 1.016064868271814
$ws.Cells.Item(22, 14).Value = 1.011999683167334
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.01748198966302
$ws.Cells.Item(23, 4).Value = 1.019829645181557
$ws.Cells.Item(23, 5).Value = 1.018962836216676
$ws.Cells.Item(23, 6).Value = 1.012984898563937
$ws.Cells.Item(23, 9).Value = 1.028433586933674
$ws.Cells.Item(23, 10).Value = 1.024321281103696
$ws.Cells.Item(23, 11).Value = 1.023532822709684
$ws.Cells.Item(23, 12).Value = 1.022669417906964
$ws.Cells.Item(23, 13).Value = 1.016715133454914
$ws.Cells.Item(23, 14).Value = 1.012114542929812
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.019362760257442
$ws.Cells.Item(24, 4).Value = 1.021476154899363
$ws.Cells.Item(24, 5).Value = 1.02056544451145
$ws.Cells.Item(24, 6).Value = 1.015826779587472
$ws.Cells.Item(24, 9).Value = 1.028817698418006
$ws.Cells.Item(24, 10).Value = 1.025666338632306
$ws.Cells.Item(24, 11).Value = 1.024897694684196
$ws.Cells.Item(24, 12).Value = 1.023990266098775
$ws.Cells.Item(24, 13).Value = 1.019268782158249
$ws.Cells.Item(24, 14).Value = 1.012564860854818
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.0215369169012
$ws.Cells.Item(25, 4).Value = 1.023380672994001
$ws.Cells.Item(25, 5).Value = 1.022419817328633
$ws.Cells.Item(25, 6).Value = 1.019113554465251
$ws.Cells.Item(25, 9).Value = 1.029250585069787
$ws.Cells.Item(25, 10).Value = 1.027216980062609
$ws.Cells.Item(25, 11).Value = 1.026472751124795
$ws.Cells.Item(25, 12).Value = 1.025514999109351
$ws.Cells.Item(25, 13).Value = 1.022219464762839
$ws.Cells.Item(25, 14).Value = 1.013083603989592
